$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: swap the STM32F746ZGT6 part for the STM32F207VGT6 part ---
$ws.Range("A2").Value = "STM32F207VGT6"
$ws.Range("F2").Value = "https://www.digikey.com/en/products/detail/stmicroelectronics/STM32F207VGT6TR/4357621"

# F10 already carried a (now pointless) "applyFill" cell style from the old
# workbook; reset it to Normal first so every hyperlinked cell in F2:F10
# ends up sharing the same new Hyperlink style.
$ws.Range("F10").Style = "Normal"

# --- Turn the plain-text Digikey URLs in column F into real hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.digikey.com/en/products/detail/stmicroelectronics/STM32F207VGT6TR/4357621")

$lm2678_5v_url = "https://www.digikey.com/en/products/detail/texas-instruments/LM2678S-5.0%2FNOPB/363825?utm_adgroup=General&utm_source=google&utm_medium=cpc&utm_campaign=PMax%20Shopping_Product_Zombie%20SKUs&utm_term=&utm_content=General&utm_id=go_cmp-17815035045_adg-_ad-__dev-c_ext-_prd-363825_sig-CjwKCAiA5L2tBhBTEiwAdSxJX2jmx8jM-JlHhV04F58rlCzi0KZgwJl8jmcjRGNCM7uSaMTsq63izRoCBfYQAvD_BwE&gad_source=1&gclid=CjwKCAiA5L2tBhBTEiwAdSxJX2jmx8jM-JlHhV04F58rlCzi0KZgwJl8jmcjRGNCM7uSaMTsq63izRoCBfYQAvD_BwE"
$ws.Hyperlinks.Add($ws.Range("F3"), $lm2678_5v_url, "", "", $lm2678_5v_url)

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.digikey.com/en/products/detail/texas-instruments/LM2678S-3-3-NOPB/366918")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.digikey.com/en/products/detail/nichicon/UUD1H150MCL1GS/590040")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.digikey.com/en/products/detail/nichicon/UWP1HR47MCL1GB/2550802")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.digikey.com/en/products/detail/vishay-general-semiconductor-diodes-division/VS-6TQ045S-M3/5426222")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.digikey.com/en/products/detail/bourns-inc/2300LL-220-V-RC/725902")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.digikey.com/en/products/detail/nichicon/UCZ1J181MNJ1MS/5144110")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.digikey.com/en/products/detail/murata-electronics/GRM2195C1H103JA01D/586788")

# --- Restore the selection to A2 (was parked at J21) ---
$ws.Range("A2").Select() | Out-Null
